# Weekly price update: a new week's record (Ajo, Feria Lagunitas de Puerto
# Montt) is inserted as the new first data row of the series (row 433),
# pushing the existing rows 433:455 down to 434:456.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 433; everything below (433:455) shifts to 434:456.
$ws.Rows.Item(433).Insert()

# Populate the new row 433 with this week's record.
$ws.Range("A433").Value = 4
$ws.Range("B433").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C433").Value = "Los Lagos"
$ws.Range("D433").Value = 45041
$ws.Range("E433").Value = 10
$ws.Range("F433").Value = 100112003
$ws.Range("G433").Value = "Ajo"
$ws.Range("H433").Value = "Chino"
$ws.Range("I433").Value = "Primera"
$ws.Range("J433").Value = 240
$ws.Range("K433").Value = 18000
$ws.Range("L433").Value = 21000
$ws.Range("M433").Value = 19500
$ws.Range("N433").Value = "$/caja 10 kilos"
$ws.Range("O433").Value = "China"
$ws.Range("P433").Value = 1950
$ws.Range("Q433").Value = 10
$ws.Range("R433").Value = "Hortaliza"
